$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 35-38 (no longer present in the updated dataset)
$ws.Range("A35:E38").Delete()

$ws.Range("A2").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B2").Value = "701号直流"
$ws.Range("C2").Value = 45927.457337962966
$ws.Range("D2").Value = 45933.341527777775
$ws.Range("E2").Value = 141.22055555542465

$ws.Range("A3").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B3").Value = "502号直流"
$ws.Range("C3").Value = 45930.238043981481
$ws.Range("D3").Value = 45933.341527777775
$ws.Range("E3").Value = 74.48361111106351

$ws.Range("A4").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B4").Value = "112号直流"
$ws.Range("C4").Value = 45930.517060185186
$ws.Range("D4").Value = 45933.341527777775
$ws.Range("E4").Value = 67.787222222134005

$ws.Range("A5").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B5").Value = "111号直流"
$ws.Range("C5").Value = 45930.618518518517
$ws.Range("D5").Value = 45933.341527777775
$ws.Range("E5").Value = 65.352222222194541

$ws.Range("A6").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B6").Value = "602号直流"
$ws.Range("C6").Value = 45931.226261574076
$ws.Range("D6").Value = 45933.341527777775
$ws.Range("E6").Value = 50.766388888761867

$ws.Range("A7").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B7").Value = "201号直流"
$ws.Range("C7").Value = 45931.575543981482
$ws.Range("D7").Value = 45933.341527777775
$ws.Range("E7").Value = 42.383611111028586

$ws.Range("A8").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B8").Value = "401号直流"
$ws.Range("C8").Value = 45931.592060185183
$ws.Range("D8").Value = 45933.341527777775
$ws.Range("E8").Value = 41.987222222203854

$ws.Range("A9").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B9").Value = "505号直流"
$ws.Range("C9").Value = 45931.647233796299
$ws.Range("D9").Value = 45933.341527777775
$ws.Range("E9").Value = 40.66305555542931

$ws.Range("A10").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B10").Value = "102号直流"
$ws.Range("C10").Value = 45931.647812499999
$ws.Range("D10").Value = 45933.341527777775
$ws.Range("E10").Value = 40.649166666611563

$ws.Range("A11").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B11").Value = "003B号直流"
$ws.Range("C11").Value = 45932.050335648149
$ws.Range("D11").Value = 45933.341527777775
$ws.Range("E11").Value = 30.988611111009959

$ws.Range("A12").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B12").Value = "801号直流"
$ws.Range("C12").Value = 45932.072800925926
$ws.Range("D12").Value = 45933.341527777775
$ws.Range("E12").Value = 30.449444444384426

$ws.Range("A13").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B13").Value = "103号直流"
$ws.Range("C13").Value = 45932.075300925928
$ws.Range("D13").Value = 45933.341527777775
$ws.Range("E13").Value = 30.389444444328547

$ws.Range("A14").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B14").Value = "603号直流"
$ws.Range("C14").Value = 45932.081099537034
$ws.Range("D14").Value = 45933.341527777775
$ws.Range("E14").Value = 30.250277777784504

$ws.Range("A15").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B15").Value = "201号直流"
$ws.Range("C15").Value = 45932.114317129628
$ws.Range("D15").Value = 45933.341527777775
$ws.Range("E15").Value = 29.45305555552477

$ws.Range("A16").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B16").Value = "B03号直流"
$ws.Range("C16").Value = 45932.143969907411
$ws.Range("D16").Value = 45933.341527777775
$ws.Range("E16").Value = 28.741388888738584

$ws.Range("A17").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B17").Value = "905号直流"
$ws.Range("C17").Value = 45932.25236111111
$ws.Range("D17").Value = 45933.341527777775
$ws.Range("E17").Value = 26.139999999955762

$ws.Range("A18").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B18").Value = "402号直流"
$ws.Range("C18").Value = 45932.255543981482
$ws.Range("D18").Value = 45933.341527777775
$ws.Range("E18").Value = 26.063611111021601

$ws.Range("A19").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B19").Value = "604号直流"
$ws.Range("C19").Value = 45932.258009259262
$ws.Range("D19").Value = 45933.341527777775
$ws.Range("E19").Value = 26.004444444319233

$ws.Range("A20").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B20").Value = "904号直流"
$ws.Range("C20").Value = 45932.463321759256
$ws.Range("D20").Value = 45933.341527777775
$ws.Range("E20").Value = 21.076944444444962

$ws.Range("A21").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B21").Value = "304号直流"
$ws.Range("C21").Value = 45932.531018518515
$ws.Range("D21").Value = 45933.341527777775
$ws.Range("E21").Value = 19.452222222229466

$ws.Range("A22").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B22").Value = "404号直流"
$ws.Range("C22").Value = 45932.532604166663
$ws.Range("D22").Value = 45933.341527777775
$ws.Range("E22").Value = 19.414166666683741

$ws.Range("A23").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B23").Value = "905号直流"
$ws.Range("C23").Value = 45932.53601851852
$ws.Range("D23").Value = 45933.341527777775
$ws.Range("E23").Value = 19.332222222117707

$ws.Range("A24").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B24").Value = "903号直流"
$ws.Range("C24").Value = 45932.538726851853
$ws.Range("D24").Value = 45933.341527777775
$ws.Range("E24").Value = 19.267222222115379

$ws.Range("A25").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B25").Value = "006A号直流"
$ws.Range("C25").Value = 45932.554108796299
$ws.Range("D25").Value = 45933.341527777775
$ws.Range("E25").Value = 18.89805555541534

$ws.Range("A26").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B26").Value = "306号直流"
$ws.Range("C26").Value = 45932.55982638889
$ws.Range("D26").Value = 45933.341527777775
$ws.Range("E26").Value = 18.76083333324641

$ws.Range("A27").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B27").Value = "210号直流"
$ws.Range("C27").Value = 45932.573414351849
$ws.Range("D27").Value = 45933.341527777775
$ws.Range("E27").Value = 18.434722222213168

$ws.Range("A28").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B28").Value = "203号直流"
$ws.Range("C28").Value = 45932.585451388892
$ws.Range("D28").Value = 45933.341527777775
$ws.Range("E28").Value = 18.145833333197515

$ws.Range("A29").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B29").Value = "804号直流"
$ws.Range("C29").Value = 45932.644085648149
$ws.Range("D29").Value = 45933.341527777775
$ws.Range("E29").Value = 16.738611111009959

$ws.Range("A30").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B30").Value = "208号直流"
$ws.Range("C30").Value = 45932.663437499999
$ws.Range("D30").Value = 45933.341527777775
$ws.Range("E30").Value = 16.274166666611563

$ws.Range("A31").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B31").Value = "901号直流"
$ws.Range("C31").Value = 45932.677557870367
$ws.Range("D31").Value = 45933.341527777775
$ws.Range("E31").Value = 15.935277777782176

$ws.Range("A32").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B32").Value = "A03号直流"
$ws.Range("C32").Value = 45932.698055555556
$ws.Range("D32").Value = 45933.341527777775
$ws.Range("E32").Value = 15.443333333241753

$ws.Range("A33").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B33").Value = "204号直流"
$ws.Range("C33").Value = 45932.723483796297
$ws.Range("D33").Value = 45933.341527777775
$ws.Range("E33").Value = 14.833055555471219

$ws.Range("A34").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B34").Value = "105号直流"
$ws.Range("C34").Value = 45932.736400462964
$ws.Range("D34").Value = 45933.341527777775
$ws.Range("E34").Value = 14.523055555473547

# Update selection to match the final saved view state
$ws.Range("G5").Select()